# Daniel and Acrolinx feedback
# Apply text/formatting fixes to the architecture diagram on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($shapes, [string]$name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

# 1) "EventBridge " + "rule"  ->  single run "EventBridge rule"
$shpEventBridge = Get-ShapeByName $s.Shapes "Google Shape;68;p13"
$trEventBridge = $shpEventBridge.TextFrame.TextRange
$trEventBridge.Text = "TEMP"
$trEventBridge.Text = "EventBridge rule"

# 2) "System Manager Automation runbook" -> "Systems " + "Manager Automation runbook"
$shpRunbook = Get-ShapeByName $s.Shapes "Google Shape;71;p13"
$trRunbook = $shpRunbook.TextFrame.TextRange
$trRunbook.Characters(1, 7).Text = "Systems "

# 3) "Systems Manager Command " + "d" + "ocument" -> single run "Systems Manager Command document"
$shpCommandDoc = Get-ShapeByName $s.Shapes "Google Shape;72;p13"
$trCommandDoc = $shpCommandDoc.TextFrame.TextRange
$trCommandDoc.Text = "TEMP"
$trCommandDoc.Text = "Systems Manager Command document"

# 4) "Systems Manager " + "Automation " + "execution " + "r" + "ole"
#    -> "Systems Manager " + "Automation execution IAM role"
$shpExecRole = Get-ShapeByName $s.Shapes "Google Shape;75;p13"
$trExecRole = $shpExecRole.TextFrame.TextRange
$trExecRole.Characters(17, 25).Text = "TEMP"
$trExecRole.Characters(17, 4).Text = "Automation execution IAM role"

# 5) "AWS IAM" shape inside Group 9 - drop the stray endParaRPr font override
$shpGroup9 = Get-ShapeByName $s.Shapes "Group 9"
$shpAwsIam = Get-ShapeByName $shpGroup9.GroupItems "TextBox 12"
$trAwsIam = $shpAwsIam.TextFrame.TextRange
$trAwsIam.Delete()
$trAwsIam.Text = "AWS IAM"

# 6) "New Relic license key" shape - drop the stray endParaRPr font override
$shpLicenseKey = Get-ShapeByName $s.Shapes "TextBox 6"
$trLicenseKey = $shpLicenseKey.TextFrame.TextRange
$trLicenseKey.Delete()
$trLicenseKey.Text = "New Relic license key"
